$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 64
$ws_ALC.Range("H64").Value = 3776.923
$ws_ALC.Range("I64").Value = 3511.111
$ws_ALC.Range("K64").Value = 3511.111
$ws_ALC.Range("M64").Value = -3263.111

# ALC row 67
$ws_ALC.Range("H67").Value = 3776.923
$ws_ALC.Range("I67").Value = 3511.111
$ws_ALC.Range("K67").Value = 3511.111
$ws_ALC.Range("M67").Value = -2653.111

# ALC row 86
$ws_ALC.Range("H86").Value = 38471736
$ws_ALC.Range("J86").Value = 83335990
$ws_ALC.Range("L86").Value = 83335990
$ws_ALC.Range("N86").Value = -83338236

# ALC row 89
$ws_ALC.Range("H89").Value = 38471736
$ws_ALC.Range("J89").Value = 83335990
$ws_ALC.Range("L89").Value = 416679950
$ws_ALC.Range("N89").Value = -416691182

# ALC row 138
$ws_ALC.Range("H138").Value = 2858.6956
$ws_ALC.Range("I138").Value = 2219.4285
$ws_ALC.Range("J138").Value = 3853.111
$ws_ALC.Range("K138").Value = 6658.2855
$ws_ALC.Range("L138").Value = 11559.333
$ws_ALC.Range("M138").Value = -1518.2855
$ws_ALC.Range("N138").Value = -21839.333

# ALC row 141
$ws_ALC.Range("H141").Value = 1567
$ws_ALC.Range("I141").Value = 1384.8214
$ws_ALC.Range("J141").Value = 2133.7778
$ws_ALC.Range("K141").Value = 4154.4642
$ws_ALC.Range("L141").Value = 6401.3334
$ws_ALC.Range("M141").Value = 1025.5358
$ws_ALC.Range("N141").Value = -16761.3334

# ARM row 55
$ws_ARM.Range("H55").Value = 16053
$ws_ARM.Range("J55").Value = 16053
$ws_ARM.Range("L55").Value = 16053
$ws_ARM.Range("N55").Value = -16683

# ARM row 132
$ws_ARM.Range("H132").Value = 94416.59
$ws_ARM.Range("I132").Value = 120581.48
$ws_ARM.Range("J132").Value = 2839.5
$ws_ARM.Range("K132").Value = 361744.44
$ws_ARM.Range("L132").Value = 8518.5
$ws_ARM.Range("M132").Value = -359214.44
$ws_ARM.Range("N132").Value = -13578.5

# CRP row 31
$ws_CRP.Range("H31").Value = 1911.8462
$ws_CRP.Range("I31").Value = 1453.0476
$ws_CRP.Range("J31").Value = 3838.8
$ws_CRP.Range("K31").Value = 1453.0476
$ws_CRP.Range("L31").Value = 3838.8
$ws_CRP.Range("M31").Value = -1158.0476
$ws_CRP.Range("N31").Value = -4428.8

# CRP row 34
$ws_CRP.Range("H34").Value = 1911.8462
$ws_CRP.Range("I34").Value = 1453.0476
$ws_CRP.Range("J34").Value = 3838.8
$ws_CRP.Range("K34").Value = 1453.0476
$ws_CRP.Range("L34").Value = 3838.8
$ws_CRP.Range("M34").Value = -1251.0476
$ws_CRP.Range("N34").Value = -4242.8

# CRP row 58
$ws_CRP.Range("H58").Value = 869.0526
$ws_CRP.Range("I58").Value = 860.5161000000001
$ws_CRP.Range("J58").Value = 906.8570999999999
$ws_CRP.Range("K58").Value = 860.5161000000001
$ws_CRP.Range("L58").Value = 906.8570999999999
$ws_CRP.Range("M58").Value = -657.5161000000001
$ws_CRP.Range("N58").Value = -1312.8571

# CRP row 62
$ws_CRP.Range("H62").Value = 2550
$ws_CRP.Range("I62").Value = 2480
$ws_CRP.Range("J62").Value = 2666.6667
$ws_CRP.Range("K62").Value = 2480
$ws_CRP.Range("L62").Value = 2666.6667
$ws_CRP.Range("M62").Value = -1856
$ws_CRP.Range("N62").Value = -3914.6667

# CRP row 65
$ws_CRP.Range("H65").Value = 2550
$ws_CRP.Range("I65").Value = 2480
$ws_CRP.Range("J65").Value = 2666.6667
$ws_CRP.Range("K65").Value = 12400
$ws_CRP.Range("L65").Value = 13333.3335
$ws_CRP.Range("M65").Value = -9280
$ws_CRP.Range("N65").Value = -19573.3335

# CRP row 99
$ws_CRP.Range("H99").Value = 2660
$ws_CRP.Range("I99").Value = 4152
$ws_CRP.Range("J99").Value = 1764.8
$ws_CRP.Range("K99").Value = 4152
$ws_CRP.Range("L99").Value = 1764.8
$ws_CRP.Range("M99").Value = -2654
$ws_CRP.Range("N99").Value = -4760.8

# CRP row 126
$ws_CRP.Range("H126").Value = 2660
$ws_CRP.Range("I126").Value = 4152
$ws_CRP.Range("J126").Value = 1764.8
$ws_CRP.Range("K126").Value = 12456
$ws_CRP.Range("L126").Value = 5294.4
$ws_CRP.Range("M126").Value = -9986
$ws_CRP.Range("N126").Value = -10234.4

# CRP row 132
$ws_CRP.Range("H132").Value = 2374.9092
$ws_CRP.Range("I132").Value = 2013.8889
$ws_CRP.Range("K132").Value = 6041.6667
$ws_CRP.Range("M132").Value = -3511.6667

# CRP row 134
$ws_CRP.Range("H134").Value = 4812.7646
$ws_CRP.Range("I134").Value = 5093.8965
$ws_CRP.Range("J134").Value = 3182.2
$ws_CRP.Range("K134").Value = 15281.6895
$ws_CRP.Range("L134").Value = 9546.599999999999
$ws_CRP.Range("M134").Value = -12746.6895
$ws_CRP.Range("N134").Value = -14616.6

# CRP row 136
$ws_CRP.Range("H136").Value = 869.0526
$ws_CRP.Range("I136").Value = 860.5161000000001
$ws_CRP.Range("J136").Value = 906.8570999999999
$ws_CRP.Range("K136").Value = 2581.5483
$ws_CRP.Range("L136").Value = 2720.5713
$ws_CRP.Range("M136").Value = -31.54830000000038
$ws_CRP.Range("N136").Value = -7820.5713

# CUL row 55
$ws_CUL.Range("H55").Value = 145143.58
$ws_CUL.Range("I55").Value = 0
$ws_CUL.Range("J55").Value = 145143.58
$ws_CUL.Range("K55").Value = 0
$ws_CUL.Range("L55").Value = 435430.74
$ws_CUL.Range("M55").ClearContents()
$ws_CUL.Range("N55").Value = -435784.74

# CUL row 69
$ws_CUL.Range("H69").Value = 2407.318
$ws_CUL.Range("I69").Value = 464
$ws_CUL.Range("J69").Value = 3517.7856
$ws_CUL.Range("K69").Value = 1392
$ws_CUL.Range("L69").Value = 10553.3568
$ws_CUL.Range("M69").Value = -581
$ws_CUL.Range("N69").Value = -12175.3568

# CUL row 72
$ws_CUL.Range("H72").Value = 2407.318
$ws_CUL.Range("I72").Value = 464
$ws_CUL.Range("J72").Value = 3517.7856
$ws_CUL.Range("K72").Value = 4176
$ws_CUL.Range("L72").Value = 31660.0704
$ws_CUL.Range("M72").Value = -120
$ws_CUL.Range("N72").Value = -39772.0704

# CUL row 80
$ws_CUL.Range("H80").Value = 5599.8335
$ws_CUL.Range("I80").Value = 3633.3333
$ws_CUL.Range("J80").Value = 5993.1333
$ws_CUL.Range("K80").Value = 10899.9999
$ws_CUL.Range("L80").Value = 17979.3999
$ws_CUL.Range("M80").Value = -9963.999899999999
$ws_CUL.Range("N80").Value = -19851.3999

# CUL row 83
$ws_CUL.Range("H83").Value = 5599.8335
$ws_CUL.Range("I83").Value = 3633.3333
$ws_CUL.Range("J83").Value = 5993.1333
$ws_CUL.Range("K83").Value = 32699.9997
$ws_CUL.Range("L83").Value = 53938.1997
$ws_CUL.Range("M83").Value = -28019.9997
$ws_CUL.Range("N83").Value = -63298.1997

# CUL row 129
$ws_CUL.Range("H129").Value = 2023.826
$ws_CUL.Range("J129").Value = 2023.238
$ws_CUL.Range("L129").Value = 6069.714
$ws_CUL.Range("N129").Value = -16069.714

# CUL row 131
$ws_CUL.Range("H131").Value = 918.36365
$ws_CUL.Range("I131").Value = 502
$ws_CUL.Range("J131").Value = 940.5106
$ws_CUL.Range("K131").Value = 1506
$ws_CUL.Range("L131").Value = 2821.5318
$ws_CUL.Range("M131").Value = 3534
$ws_CUL.Range("N131").Value = -12901.5318

# GSM row 132
$ws_GSM.Range("H132").Value = 2132.5518
$ws_GSM.Range("I132").Value = 1808.619
$ws_GSM.Range("K132").Value = 5425.857
$ws_GSM.Range("M132").Value = -2895.857

# LTW row 62
$ws_LTW.Range("H62").Value = 28124.5
$ws_LTW.Range("J62").Value = 28124.5
$ws_LTW.Range("L62").Value = 28124.5
$ws_LTW.Range("N62").Value = -29372.5

# LTW row 65
$ws_LTW.Range("H65").Value = 28124.5
$ws_LTW.Range("J65").Value = 28124.5
$ws_LTW.Range("L65").Value = 84373.5
$ws_LTW.Range("N65").Value = -90613.5

# LTW row 132
$ws_LTW.Range("H132").Value = 1550
$ws_LTW.Range("I132").Value = 1203.0667
$ws_LTW.Range("J132").Value = 2851
$ws_LTW.Range("K132").Value = 3609.2001
$ws_LTW.Range("L132").Value = 8553
$ws_LTW.Range("M132").Value = -1079.2001
$ws_LTW.Range("N132").Value = -13613

# WVR row 40
$ws_WVR.Range("H40").Value = 15196
$ws_WVR.Range("J40").Value = 15196
$ws_WVR.Range("L40").Value = 15196
$ws_WVR.Range("N40").Value = -15494

# WVR row 70
$ws_WVR.Range("H70").Value = 12700
$ws_WVR.Range("I70").Value = 12700
$ws_WVR.Range("K70").Value = 12700
$ws_WVR.Range("M70").Value = -12385

# WVR row 73
$ws_WVR.Range("H73").Value = 12700
$ws_WVR.Range("I73").Value = 12700
$ws_WVR.Range("K73").Value = 12700
$ws_WVR.Range("M73").Value = -11608

# WVR row 132
$ws_WVR.Range("H132").Value = 1968.875
$ws_WVR.Range("I132").Value = 1816.3823
$ws_WVR.Range("K132").Value = 5449.1469
$ws_WVR.Range("M132").Value = -2919.1469
